# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to H:N columns across several sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 315.3
$ws.Range("I53").Value = 155.6
$ws.Range("K53").Value = 155.6
$ws.Range("M53").Value = 481.4
$ws.Range("H87").Value = 54328.5
$ws.Range("J87").Value = 54328.5
$ws.Range("L87").Value = 54328.5
$ws.Range("N87").Value = -56824.5
$ws.Range("H90").Value = 54328.5
$ws.Range("J90").Value = 54328.5
$ws.Range("L90").Value = 162985.5
$ws.Range("N90").Value = -175465.5
$ws.Range("H103").Value = 1442.0834
$ws.Range("I103").Value = 1098.75
$ws.Range("J103").Value = 1613.75
$ws.Range("K103").Value = 3296.25
$ws.Range("L103").Value = 4841.25
$ws.Range("M103").Value = -2710.25
$ws.Range("N103").Value = -6013.25
$ws.Range("H118").Value = 804.4
$ws.Range("I118").Value = 505.5
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 1516.5
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = 140.5
$ws.Range("N118").Value = -9314
$ws.Range("H127").Value = 1038
$ws.Range("I127").Value = 448.5
$ws.Range("J127").Value = 2217
$ws.Range("K127").Value = 1345.5
$ws.Range("L127").Value = 6651
$ws.Range("M127").Value = 3614.5
$ws.Range("N127").Value = -16571
$ws.Range("H138").Value = 4393.5
$ws.Range("I138").Value = 2316
$ws.Range("J138").Value = 5038.241
$ws.Range("K138").Value = 6948
$ws.Range("L138").Value = 15114.723
$ws.Range("M138").Value = -1808
$ws.Range("N138").Value = -25394.723

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1929.4615
$ws.Range("I2").Value = 1953.7778
$ws.Range("J2").Value = 1874.75
$ws.Range("K2").Value = 1953.7778
$ws.Range("L2").Value = 1874.75
$ws.Range("M2").Value = -1840.7778
$ws.Range("N2").Value = -2100.75
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H116").Value = 1929.4615
$ws.Range("I116").Value = 1953.7778
$ws.Range("J116").Value = 1874.75
$ws.Range("K116").Value = 1953.7778
$ws.Range("L116").Value = 1874.75
$ws.Range("M116").Value = 340.2221999999999
$ws.Range("N116").Value = -6462.75
$ws.Range("H128").Value = 96666
$ws.Range("J128").Value = 96666
$ws.Range("L128").Value = 96666
$ws.Range("N128").Value = -106626
$ws.Range("H132").Value = 2860.0667
$ws.Range("I132").Value = 2325.3333
$ws.Range("K132").Value = 6975.999899999999
$ws.Range("M132").Value = -4445.999899999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1929.4615
$ws.Range("I3").Value = 1953.7778
$ws.Range("J3").Value = 1874.75
$ws.Range("K3").Value = 1953.7778
$ws.Range("L3").Value = 1874.75
$ws.Range("M3").Value = -1839.7778
$ws.Range("N3").Value = -2102.75
$ws.Range("H38").Value = 50000
$ws.Range("J38").Value = 50000
$ws.Range("L38").Value = 50000
$ws.Range("N38").Value = -50832
$ws.Range("H62").Value = 90000
$ws.Range("J62").Value = 90000
$ws.Range("L62").Value = 90000
$ws.Range("N62").Value = -91372
$ws.Range("H65").Value = 90000
$ws.Range("J65").Value = 90000
$ws.Range("L65").Value = 270000
$ws.Range("N65").Value = -276864
$ws.Range("H105").Value = 2884.7778
$ws.Range("I105").Value = 2895.5
$ws.Range("J105").Value = 2799
$ws.Range("K105").Value = 2895.5
$ws.Range("L105").Value = 2799
$ws.Range("M105").Value = -1148.5
$ws.Range("N105").Value = -6293

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 14003
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 14003
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 14003
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -14343
$ws.Range("H31").Value = 1694.3
$ws.Range("I31").Value = 1184.7142
$ws.Range("K31").Value = 1184.7142
$ws.Range("M31").Value = -889.7141999999999
$ws.Range("H34").Value = 1694.3
$ws.Range("I34").Value = 1184.7142
$ws.Range("K34").Value = 1184.7142
$ws.Range("M34").Value = -982.7141999999999
$ws.Range("H35").Value = 26250
$ws.Range("I35").Value = 2500
$ws.Range("K35").Value = 2500
$ws.Range("M35").Value = -2206
$ws.Range("H50").Value = 38281.668
$ws.Range("J50").Value = 38281.668
$ws.Range("L50").Value = 38281.668
$ws.Range("N50").Value = -39531.668
$ws.Range("H60").Value = 37586
$ws.Range("J60").Value = 46780.332
$ws.Range("L60").Value = 46780.332
$ws.Range("N60").Value = -47802.332
$ws.Range("H86").Value = 4936.7646
$ws.Range("J86").Value = 4312.0835
$ws.Range("L86").Value = 4312.0835
$ws.Range("N86").Value = -6558.0835
$ws.Range("H89").Value = 4936.7646
$ws.Range("J89").Value = 4312.0835
$ws.Range("L89").Value = 21560.4175
$ws.Range("N89").Value = -32792.4175
$ws.Range("H134").Value = 2751.5
$ws.Range("I134").Value = 2662.25
$ws.Range("K134").Value = 7986.75
$ws.Range("M134").Value = -5451.75
$ws.Range("H141").Value = 236625
$ws.Range("J141").Value = 236625
$ws.Range("L141").Value = 236625
$ws.Range("N141").Value = -246985

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2499.75
$ws.Range("J114").Value = 2333.3333
$ws.Range("L114").Value = 6999.999899999999
$ws.Range("N114").Value = -13507.9999
$ws.Range("H137").Value = 2975.1333
$ws.Range("J137").Value = 3458.6667
$ws.Range("L137").Value = 10376.0001
$ws.Range("N137").Value = -20576.0001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3811.2144
$ws.Range("I132").Value = 3237.8572
$ws.Range("J132").Value = 4384.5713
$ws.Range("K132").Value = 9713.571599999999
$ws.Range("L132").Value = 13153.7139
$ws.Range("M132").Value = -7183.571599999999
$ws.Range("N132").Value = -18213.7139

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 650
$ws.Range("J55").Value = 597.8570999999999
$ws.Range("L55").Value = 597.8570999999999
$ws.Range("N55").Value = -943.8570999999999
$ws.Range("H132").Value = 3699.6667
$ws.Range("I132").Value = 2550
$ws.Range("K132").Value = 7650
$ws.Range("M132").Value = -5120

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3475.3333
$ws.Range("I132").Value = 2713.5
$ws.Range("K132").Value = 8140.5
$ws.Range("M132").Value = -5610.5
